$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8 & 9: Count/Total functions: Int32 -> ExpressionResultNumeric ---
$ws.Range("B8").Value = "ExpressionResultNumeric"
$ws.Range("B9").Value = "ExpressionResultNumeric"

# --- Row 43: Pluralize() loses its Int32 overload (remove G43) ---
$ws.Range("G43").Clear()

# --- Rows 77-83: Int32 -> ExpressionResultNumeric ---
$ws.Range("B77").Value = "ExpressionResultNumeric"
$ws.Range("B78").Value = "ExpressionResultNumeric"
$ws.Range("B79").Value = "ExpressionResultNumeric"
$ws.Range("B80").Value = "ExpressionResultNumeric"
$ws.Range("B81").Value = "ExpressionResultNumeric"
$ws.Range("B82").Value = "ExpressionResultNumeric"
$ws.Range("B83").Value = "ExpressionResultNumeric"

# --- Row 143-146: Int32 -> ExpressionResultNumeric ---
$ws.Range("D143").Value = "ExpressionResultNumeric"
$ws.Range("B144").Value = "ExpressionResultNumeric"
$ws.Range("B145").Value = "ExpressionResultNumeric"
$ws.Range("B146").Value = "ExpressionResultNumeric"

# --- New rows 152-154: Year / Month / Day DateTime functions ---
$names = @("Year", "Month", "Day")
foreach ($i in 0..2) {
  $r = 152 + $i

  $ws.Range("A146:P146").Copy()
  $ws.Range("A$r").PasteSpecial(-4122)
  $ws.Range("C$r").Clear()
  $ws.Range("E$r").Clear()
  $ws.Range("F$r").Clear()
  $ws.Range("G$r").Clear()
  $ws.Range("H$r").Clear()

  $ws.Range("A$r").Value = $names[$i]
  $ws.Range("B$r").Value = "ExpressionResultNumeric"
  $ws.Range("D$r").Value = "DateTime"

  $ws.Range("I$r").Formula = '=CONCATENATE("""",A' + $r + ',"""")'
  $ws.Range("J$r").Formula = '=CONCATENATE("""",B' + $r + ',"""")'
  $ws.Range("K$r").Formula = '=CONCATENATE("""",D' + $r + ',"""")'
  $ws.Range("L$r").Formula = '=CONCATENATE("""",E' + $r + ',"""")'
  $ws.Range("M$r").Formula = '=CONCATENATE("""",F' + $r + ',"""")'
  $ws.Range("N$r").Formula = '=CONCATENATE("""",G' + $r + ',"""")'
  $ws.Range("O$r").Formula = '=CONCATENATE("""",H' + $r + ',"""")'
  $ws.Range("P$r").Formula = '=SUBSTITUTE(CONCATENATE($Q$1,I' + $r + ',$P$1,J' + $r + ',$P$1,K' + $r + ',$P$1,L' + $r + ',$P$1,M' + $r + ',$P$1,N' + $r + ',$P$1,O' + $r + ',$R$1),", """"","")'
}

# --- AutoFilter range grows to include the new rows ---
$ws.AutoFilterMode = $false
$ws.Range("A1:H154").AutoFilter()

# --- _FilterDatabase defined name grows to match ---
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$H`$154"

# --- Selection moved to C11 ---
$ws.Range("C11").Select()
